{"js": "// Update the date line and the 25 multiplication answers in the table\n// to the new values, matching the target revision of the worksheet.\nconst replacements = [\n  [\"2025-03-19 Wednesday\", \"2025-03-20 Thursday\"],\n  [\"77\u00d714=1078\", \"53\u00d767=3551\"],\n  [\"47\u00d712=564\", \"17\u00d794=1598\"],\n  [\"81\u00d789=7209\", \"25\u00d730=750\"],\n  [\"99\u00d741=4059\", \"59\u00d785=5015\"],\n  [\"40\u00d769=2760\", \"69\u00d728=1932\"],\n  [\"47\u00d799=4653\", \"13\u00d742=546\"],\n  [\"87\u00d722=1914\", \"19\u00d724=456\"],\n  [\"60\u00d765=3900\", \"89\u00d767=5963\"],\n  [\"92\u00d748=4416\", \"95\u00d750=4750\"],\n  [\"64\u00d748=3072\", \"72\u00d719=1368\"],\n  [\"50\u00d745=2250\", \"14\u00d729=406\"],\n  [\"89\u00d765=5785\", \"22\u00d747=1034\"],\n  [\"84\u00d759=4956\", \"46\u00d771=3266\"],\n  [\"60\u00d757=3420\", \"49\u00d778=3822\"],\n  [\"65\u00d717=1105\", \"93\u00d764=5952\"],\n  [\"68\u00d721=1428\", \"20\u00d787=1740\"],\n  [\"60\u00d718=1080\", \"53\u00d730=1590\"],\n  [\"87\u00d716=1392\", \"16\u00d777=1232\"],\n  [\"81\u00d788=7128\", \"91\u00d760=5460\"],\n  [\"89\u00d714=1246\", \"30\u00d776=2280\"],\n  [\"63\u00d744=2772\", \"11\u00d748=528\"],\n  [\"77\u00d712=924\", \"11\u00d757=627\"],\n  [\"31\u00d751=1581\", \"78\u00d730=2340\"],\n  [\"35\u00d717=595\", \"71\u00d762=4402\"],\n  [\"12\u00d788=1056\", \"68\u00d758=3944\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication answers in the table\n# to the new values, matching the target revision of the worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-19 Wednesday\", \"2025-03-20 Thursday\"),\n    @(\"77\u00d714=1078\", \"53\u00d767=3551\"),\n    @(\"47\u00d712=564\", \"17\u00d794=1598\"),\n    @(\"81\u00d789=7209\", \"25\u00d730=750\"),\n    @(\"99\u00d741=4059\", \"59\u00d785=5015\"),\n    @(\"40\u00d769=2760\", \"69\u00d728=1932\"),\n    @(\"47\u00d799=4653\", \"13\u00d742=546\"),\n    @(\"87\u00d722=1914\", \"19\u00d724=456\"),\n    @(\"60\u00d765=3900\", \"89\u00d767=5963\"),\n    @(\"92\u00d748=4416\", \"95\u00d750=4750\"),\n    @(\"64\u00d748=3072\", \"72\u00d719=1368\"),\n    @(\"50\u00d745=2250\", \"14\u00d729=406\"),\n    @(\"89\u00d765=5785\", \"22\u00d747=1034\"),\n    @(\"84\u00d759=4956\", \"46\u00d771=3266\"),\n    @(\"60\u00d757=3420\", \"49\u00d778=3822\"),\n    @(\"65\u00d717=1105\", \"93\u00d764=5952\"),\n    @(\"68\u00d721=1428\", \"20\u00d787=1740\"),\n    @(\"60\u00d718=1080\", \"53\u00d730=1590\"),\n    @(\"87\u00d716=1392\", \"16\u00d777=1232\"),\n    @(\"81\u00d788=7128\", \"91\u00d760=5460\"),\n    @(\"89\u00d714=1246\", \"30\u00d776=2280\"),\n    @(\"63\u00d744=2772\", \"11\u00d748=528\"),\n    @(\"77\u00d712=924\", \"11\u00d757=627\"),\n    @(\"31\u00d751=1581\", \"78\u00d730=2340\"),\n    @(\"35\u00d717=595\", \"71\u00d762=4402\"),\n    @(\"12\u00d788=1056\", \"68\u00d758=3944\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
